## Edit: Tue, Apr 14, 2020 12:05:08 AM
##
## 1. The table on the "PLENARY" slide currently uses the custom
##    table style {141DF45D-04F4-439D-84E6-BB67EBA34F20} ("Table_0",
##    defined in ppt/tableStyles.xml). Re-style it with the built-in
##    PowerPoint table style {DC8088E9-F2C7-46F7-86A1-42CDB174E751}.
##
## 2. The deck's Design theme ("Integral") and the Notes Master theme
##    ("Office Theme") are swapped, so the main slides pick up the
##    default Office look and the notes pages pick up Integral.

$p = $ppt.ActivePresentation

# --- 1. Re-apply the built-in table style on every table shape that
#        still carries the old custom style id. -----------------------
$oldStyleId = "{141DF45D-04F4-439D-84E6-BB67EBA34F20}"
$newStyleId = "{DC8088E9-F2C7-46F7-86A1-42CDB174E751}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $tbl = $shape.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}

# --- 2. Swap the Design theme applied to the slide master with the
#        one applied to the notes master: the main deck goes back to
#        the default "Office Theme" and the notes pages pick up the
#        "Integral" theme that used to drive the slides. ---------------
$slideMaster = $p.SlideMaster
$notesMaster = $p.NotesMaster

$officeThemeFile = "$env:ProgramFiles\Microsoft Office\root\Document Themes 16\Office Theme.thmx"
$integralThemeFile = "$env:ProgramFiles\Microsoft Office\root\Document Themes 16\Integral.thmx"

$slideMaster.ApplyTheme($officeThemeFile)
$notesMaster.ApplyTheme($integralThemeFile)
